{"js": "const replacements = [\n  [\"2024-08-19 Monday\", \"2024-08-20 Tuesday\"],\n  [\"93-48=\", \"46-44=\"],\n  [\"5+38=\", \"63-58=\"],\n  [\"35-21=\", \"47+25=\"],\n  [\"83-62=\", \"69+29=\"],\n  [\"0+39=\", \"61-51=\"],\n  [\"42+44=\", \"96-30=\"],\n  [\"97-62=\", \"58-23=\"],\n  [\"26-10=\", \"45-42=\"],\n  [\"87-42=\", \"81-9=\"],\n  [\"58-21=\", \"70-16=\"],\n  [\"18-2=\", \"38-13=\"],\n  [\"61-36=\", \"57-18=\"],\n  [\"22+30=\", \"66-47=\"],\n  [\"31+4=\", \"6+5=\"],\n  [\"60+38=\", \"57-15=\"],\n  [\"11+41=\", \"40+0=\"],\n  [\"91-46=\", \"62-34=\"],\n  [\"12-6=\", \"37-29=\"],\n  [\"64+9=\", \"63+22=\"],\n  [\"68-13=\", \"23+68=\"],\n  [\"6+41=\", \"17+19=\"],\n  [\"69+3=\", \"81-6=\"],\n  [\"85-57=\", \"81-28=\"],\n  [\"38-1=\", \"60+0=\"],\n  [\"60-19=\", \"10+27=\"],\n  [\"76-64=\", \"87+7=\"],\n  [\"30+19=\", \"88-31=\"],\n  [\"40+3=\", \"69-23=\"],\n  [\"23-7=\", \"11+6=\"],\n  [\"75+21=\", \"3+3=\"],\n  [\"87-55=\", \"98-4=\"],\n  [\"9+85=\", \"47-2=\"],\n  [\"42+22=\", \"54+19=\"],\n  [\"97-74=\", \"3-0=\"],\n  [\"59+16=\", \"69-2=\"],\n  [\"29+3=\", \"8+1=\"],\n  [\"94-50=\", \"31-4=\"],\n  [\"56-24=\", \"43+10=\"],\n  [\"25-1=\", \"86-76=\"],\n  [\"74-51=\", \"45-2=\"],\n  [\"2+35=\", \"65-7=\"],\n  [\"24+17=\", \"67+11=\"],\n  [\"23+62=\", \"78-77=\"],\n  [\"33+7=\", \"80-51=\"],\n  [\"74+16=\", \"4+91=\"],\n  [\"28+59=\", \"32-16=\"],\n  [\"89-30=\", \"57-42=\"],\n  [\"67+0=\", \"1+34=\"],\n  [\"21+7=\", \"60-29=\"],\n  [\"35-34=\", \"26+43=\"],\n  [\"16+80=\", \"29+40=\"],\n  [\"5+53=\", \"52+39=\"],\n  [\"25-9=\", \"75-62=\"],\n  [\"57+38=\", \"93-35=\"],\n  [\"47+31=\", \"40+46=\"],\n  [\"83-54=\", \"4+39=\"],\n  [\"8+55=\", \"31-1=\"],\n  [\"44-27=\", \"80-25=\"],\n  [\"27+72=\", \"92+0=\"],\n  [\"77-4=\", \"1+50=\"],\n  [\"71-11=\", \"33-31=\"],\n  [\"81+12=\", \"78-52=\"],\n  [\"88-38=\", \"54+14=\"],\n  [\"46-35=\", \"39+54=\"],\n  [\"4+22=\", \"80+5=\"],\n  [\"4+54=\", \"98-4=\"],\n  [\"41+19=\", \"95-21=\"],\n  [\"97-80=\", \"74+0=\"],\n  [\"27-13=\", \"12+31=\"],\n  [\"50-11=\", \"59-28=\"],\n  [\"87-39=\", \"60+24=\"],\n  [\"58-5=\", \"34+12=\"],\n  [\"28+25=\", \"30-3=\"],\n  [\"98-77=\", \"77-61=\"],\n  [\"31-17=\", \"82-75=\"],\n  [\"88-83=\", \"81+11=\"],\n  [\"61-20=\", \"70+2=\"],\n  [\"6+81=\", \"77-19=\"],\n  [\"59-32=\", \"1+20=\"],\n  [\"28+41=\", \"32+63=\"],\n  [\"55+16=\", \"58-57=\"],\n  [\"47-40=\", \"93-79=\"],\n  [\"4+55=\", \"24-13=\"],\n  [\"16+48=\", \"21+72=\"],\n  [\"13-5=\", \"7+20=\"],\n  [\"91-75=\", \"7+3=\"],\n  [\"74-15=\", \"91+2=\"],\n  [\"4+51=\", \"56-41=\"],\n  [\"82+14=\", \"2+2=\"],\n  [\"41+12=\", \"45-3=\"],\n  [\"9+32=\", \"44-36=\"],\n  [\"3+76=\", \"81-3=\"],\n  [\"49-18=\", \"71-38=\"],\n  [\"14+9=\", \"74-10=\"],\n  [\"44-11=\", \"2+61=\"],\n  [\"11+44=\", \"7+2=\"],\n  [\"46+7=\", \"15+45=\"],\n  [\"54-44=\", \"61-58=\"],\n  [\"61+8=\", \"87-5=\"],\n  [\"52-47=\", \"95-91=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  // Each target string is unique in this document, so the first (and only) hit is the right one.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$replacements = @(\n  @(\"2024-08-19 Monday\", \"2024-08-20 Tuesday\"),\n  @(\"93-48=\", \"46-44=\"),\n  @(\"5+38=\", \"63-58=\"),\n  @(\"35-21=\", \"47+25=\"),\n  @(\"83-62=\", \"69+29=\"),\n  @(\"0+39=\", \"61-51=\"),\n  @(\"42+44=\", \"96-30=\"),\n  @(\"97-62=\", \"58-23=\"),\n  @(\"26-10=\", \"45-42=\"),\n  @(\"87-42=\", \"81-9=\"),\n  @(\"58-21=\", \"70-16=\"),\n  @(\"18-2=\", \"38-13=\"),\n  @(\"61-36=\", \"57-18=\"),\n  @(\"22+30=\", \"66-47=\"),\n  @(\"31+4=\", \"6+5=\"),\n  @(\"60+38=\", \"57-15=\"),\n  @(\"11+41=\", \"40+0=\"),\n  @(\"91-46=\", \"62-34=\"),\n  @(\"12-6=\", \"37-29=\"),\n  @(\"64+9=\", \"63+22=\"),\n  @(\"68-13=\", \"23+68=\"),\n  @(\"6+41=\", \"17+19=\"),\n  @(\"69+3=\", \"81-6=\"),\n  @(\"85-57=\", \"81-28=\"),\n  @(\"38-1=\", \"60+0=\"),\n  @(\"60-19=\", \"10+27=\"),\n  @(\"76-64=\", \"87+7=\"),\n  @(\"30+19=\", \"88-31=\"),\n  @(\"40+3=\", \"69-23=\"),\n  @(\"23-7=\", \"11+6=\"),\n  @(\"75+21=\", \"3+3=\"),\n  @(\"87-55=\", \"98-4=\"),\n  @(\"9+85=\", \"47-2=\"),\n  @(\"42+22=\", \"54+19=\"),\n  @(\"97-74=\", \"3-0=\"),\n  @(\"59+16=\", \"69-2=\"),\n  @(\"29+3=\", \"8+1=\"),\n  @(\"94-50=\", \"31-4=\"),\n  @(\"56-24=\", \"43+10=\"),\n  @(\"25-1=\", \"86-76=\"),\n  @(\"74-51=\", \"45-2=\"),\n  @(\"2+35=\", \"65-7=\"),\n  @(\"24+17=\", \"67+11=\"),\n  @(\"23+62=\", \"78-77=\"),\n  @(\"33+7=\", \"80-51=\"),\n  @(\"74+16=\", \"4+91=\"),\n  @(\"28+59=\", \"32-16=\"),\n  @(\"89-30=\", \"57-42=\"),\n  @(\"67+0=\", \"1+34=\"),\n  @(\"21+7=\", \"60-29=\"),\n  @(\"35-34=\", \"26+43=\"),\n  @(\"16+80=\", \"29+40=\"),\n  @(\"5+53=\", \"52+39=\"),\n  @(\"25-9=\", \"75-62=\"),\n  @(\"57+38=\", \"93-35=\"),\n  @(\"47+31=\", \"40+46=\"),\n  @(\"83-54=\", \"4+39=\"),\n  @(\"8+55=\", \"31-1=\"),\n  @(\"44-27=\", \"80-25=\"),\n  @(\"27+72=\", \"92+0=\"),\n  @(\"77-4=\", \"1+50=\"),\n  @(\"71-11=\", \"33-31=\"),\n  @(\"81+12=\", \"78-52=\"),\n  @(\"88-38=\", \"54+14=\"),\n  @(\"46-35=\", \"39+54=\"),\n  @(\"4+22=\", \"80+5=\"),\n  @(\"4+54=\", \"98-4=\"),\n  @(\"41+19=\", \"95-21=\"),\n  @(\"97-80=\", \"74+0=\"),\n  @(\"27-13=\", \"12+31=\"),\n  @(\"50-11=\", \"59-28=\"),\n  @(\"87-39=\", \"60+24=\"),\n  @(\"58-5=\", \"34+12=\"),\n  @(\"28+25=\", \"30-3=\"),\n  @(\"98-77=\", \"77-61=\"),\n  @(\"31-17=\", \"82-75=\"),\n  @(\"88-83=\", \"81+11=\"),\n  @(\"61-20=\", \"70+2=\"),\n  @(\"6+81=\", \"77-19=\"),\n  @(\"59-32=\", \"1+20=\"),\n  @(\"28+41=\", \"32+63=\"),\n  @(\"55+16=\", \"58-57=\"),\n  @(\"47-40=\", \"93-79=\"),\n  @(\"4+55=\", \"24-13=\"),\n  @(\"16+48=\", \"21+72=\"),\n  @(\"13-5=\", \"7+20=\"),\n  @(\"91-75=\", \"7+3=\"),\n  @(\"74-15=\", \"91+2=\"),\n  @(\"4+51=\", \"56-41=\"),\n  @(\"82+14=\", \"2+2=\"),\n  @(\"41+12=\", \"45-3=\"),\n  @(\"9+32=\", \"44-36=\"),\n  @(\"3+76=\", \"81-3=\"),\n  @(\"49-18=\", \"71-38=\"),\n  @(\"14+9=\", \"74-10=\"),\n  @(\"44-11=\", \"2+61=\"),\n  @(\"11+44=\", \"7+2=\"),\n  @(\"46+7=\", \"15+45=\"),\n  @(\"54-44=\", \"61-58=\"),\n  @(\"61+8=\", \"87-5=\"),\n  @(\"52-47=\", \"95-91=\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n  if (-not $ok) {\n    throw \"Replace failed for: $oldText\"\n  }\n}"}
